# Apply data updates to the active worksheet of the workbook.
# Column layout: A=Processo, B=Cargo, C=Campus, D=Nivel,
#                E=Inscritos, F=Pagos, G=Isencoes deferidas, H=Inscricoes homologadas

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> hashtable of column letter -> new value
$updates = @{
    5  = @{ E = 136 }
    7  = @{ E = 33;  F = 17;  H = 22 }
    10 = @{ E = 584; F = 287; H = 383 }
    11 = @{ E = 368; F = 198; H = 262 }
    12 = @{ E = 568 }
    13 = @{ E = 141 }
    14 = @{ E = 127 }
    16 = @{ E = 208; F = 105; H = 153 }
    20 = @{ E = 88 }
    25 = @{ E = 277 }
    27 = @{ E = 337; F = 174; H = 255 }
    30 = @{ E = 216 }
    32 = @{ E = 186; F = 110; H = 148 }
    33 = @{ E = 302; F = 155; H = 244 }
    34 = @{ E = 221; F = 146; H = 185 }
    36 = @{ E = 76 }
    37 = @{ E = 166 }
    42 = @{ E = 389 }
    43 = @{ E = 121 }
    45 = @{ E = 149 }
    46 = @{ E = 331 }
    47 = @{ E = 463 }
    48 = @{ E = 218 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
